# Refined metadata to be additional tab
# 1) Update the "time_taken" (column F) timestamps on the "data" sheet.
# 2) Add a new "metadata" sheet (after "data") describing the panelapp
#    query that produced the data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# --- 1) Refresh F2:F23 timestamps -----------------------------------------
$ws.Range("F2").Value  = "2021-10-05 14:33:43.215421"
$ws.Range("F3").Value  = "2021-10-05 14:33:43.215429"
$ws.Range("F4").Value  = "2021-10-05 14:33:43.215432"
$ws.Range("F5").Value  = "2021-10-05 14:33:43.215435"
$ws.Range("F6").Value  = "2021-10-05 14:33:43.215437"
$ws.Range("F7").Value  = "2021-10-05 14:33:43.215440"
$ws.Range("F8").Value  = "2021-10-05 14:33:43.215443"
$ws.Range("F9").Value  = "2021-10-05 14:33:43.215445"
$ws.Range("F10").Value = "2021-10-05 14:33:43.215448"
$ws.Range("F11").Value = "2021-10-05 14:33:43.215450"
$ws.Range("F12").Value = "2021-10-05 14:33:43.215453"
$ws.Range("F13").Value = "2021-10-05 14:33:43.215455"
$ws.Range("F14").Value = "2021-10-05 14:33:43.215458"
$ws.Range("F15").Value = "2021-10-05 14:33:43.215460"
$ws.Range("F16").Value = "2021-10-05 14:33:43.215463"
$ws.Range("F17").Value = "2021-10-05 14:33:43.215466"
$ws.Range("F18").Value = "2021-10-05 14:33:43.215468"
$ws.Range("F19").Value = "2021-10-05 14:33:43.215471"
$ws.Range("F20").Value = "2021-10-05 14:33:43.215474"
$ws.Range("F21").Value = "2021-10-05 14:33:43.215476"
$ws.Range("F22").Value = "2021-10-05 14:33:43.215479"
$ws.Range("F23").Value = "2021-10-05 14:33:43.215481"

# --- 2) Add the "metadata" sheet, placed right after "data" ---------------
$meta = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$meta.Name = "metadata"

# Header row (B1:G1) reuses the same bold/bordered/centered style as the
# "data" sheet's own header row (style index carried via copy/paste of
# formats, so no duplicate style gets created).
$ws.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)

$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Row 2: A2 reuses the same style as "data"!A2 (bold/bordered/centered).
$ws.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)
$meta.Range("A2").Value = 0

$meta.Range("B2").Value = "Diamond Blackfan anaemia"
$meta.Range("C2").Value = 98

# "1.0" must stay a text value (not be coerced to the number 1) and must
# not pick up a lingering custom number format/style on the cell itself.
# Stage it via a throwaway cell formatted as text, copy just the *value*
# (xlPasteValues) into the target, then remove the helper cell entirely.
$helper = $meta.Range("Z100")
$helper.NumberFormat = "@"
$helper.Value = "1.0"
$helper.Copy()
$meta.Range("D2").PasteSpecial(-4163)
$helper.Delete()

$meta.Range("E2").Value = "2021-03-06T03:02:17.262010Z"
$meta.Range("F2").Value = "2021-10-05 14:33:43.211578"
$meta.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/98/?format=json"

[void]$ws.Select()
[void]$ws.Range("A1").Select()

Write-Output "edit applied"
